$d = $word.ActiveDocument

$replacements = @(
    @{old = "41÷8="; new = "47÷4="},
    @{old = "97÷7="; new = "56÷2="},
    @{old = "37÷2="; new = "43÷9="},
    @{old = "49÷5="; new = "92÷8="},
    @{old = "15÷4="; new = "89÷4="},
    @{old = "55÷4="; new = "58÷6="},
    @{old = "50÷7="; new = "34÷5="},
    @{old = "48÷2="; new = "65÷2="},
    @{old = "53÷4="; new = "34÷7="},
    @{old = "55÷8="; new = "64÷5="},
    @{old = "33÷3="; new = "26÷5="},
    @{old = "93÷3="; new = "37÷8="},
    @{old = "60÷9="; new = "13÷7="},
    @{old = "42÷2="; new = "86÷3="},
    @{old = "11÷6="; new = "93÷7="},
    @{old = "80÷5="; new = "77÷3="},
    @{old = "16÷2="; new = "71÷2="},
    @{old = "63÷8="; new = "32÷6="},
    @{old = "90÷3="; new = "83÷5="},
    @{old = "37÷7="; new = "46÷8="},
    @{old = "90÷8="; new = "10÷9="},
    @{old = "26÷8="; new = "73÷2="},
    @{old = "92÷6="; new = "48÷8="},
    @{old = "75÷8="; new = "84÷8="},
    @{old = "78÷3="; new = "45÷9="}
)

foreach ($pair in $replacements) {
    $range = $d.Content
    $range.Find.Execute($pair.old, $true, $false, $false, $false, $false,
                         $true, 1, $false, $pair.new, 2)
}
